$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New loading-percent values for rows 2-25 (A1-index 0-23), columns B:O
# (columns E, I, L remain 0 and are included for completeness)
$data = @(
    @(9.028827612820024, 5.175968639971936, 11.05644877430328, 0, 31.06654833518582, 30.01926160710828, 14.70340521749531, 0, 11.42470144892342, 8.317913269145972, 0, 15.66587606172371, 19.00013952161483, 22.53067573583771),
    @(8.7639846002581, 5.078149553523112, 11.01323591665296, 0, 31.12002960576156, 30.09123711214846, 14.74589491319978, 0, 11.4425112782153, 8.140788427361215, 0, 15.58838201403844, 19.05455838189547, 22.59977026413126),
    @(8.598603619572538, 5.016690329056214, 10.98862081234335, 0, 31.15949652280302, 30.14352287993433, 14.7739830155486, 0, 11.45507590535954, 8.031118459560282, 0, 15.54312141005612, 19.08957656553372, 22.64628819388755),
    @(8.53061922770293, 4.99131529361111, 10.97907995758214, 0, 31.17724481785702, 30.16685896161711, 14.7859321783121, 0, 11.46060595648801, 7.986257844927544, 0, 15.52527529551741, 19.10425141977456, 22.66627308609723),
    @(8.519298032376808, 4.987082512485241, 10.97752551662808, 0, 31.1802924232126, 30.17085628039779, 14.78794671861331, 0, 11.46154897570815, 7.978800369268139, 0, 15.52234847946458, 19.10671264173223, 22.66965365104755),
    @(8.597689004117131, 5.016349418990328, 10.98849014720252, 0, 31.15972914280809, 30.14382939032049, 14.77414212897154, 0, 11.45514882586538, 8.030514056982188, 0, 15.54287829171158, 19.08977283580327, 22.64655355418589),
    @(8.938144881908833, 5.142543431195303, 11.04115577029033, 0, 31.08361185114858, 30.04239589470331, 14.7176409020703, 0, 11.43050422021911, 8.257066255491601, 0, 15.63868318170676, 19.01857075123857, 22.55364936523258),
    @(9.579566762994093, 5.378017785762387, 11.15928430406101, 0, 30.98700892211799, 29.90792645051743, 14.62269174981832, 0, 11.39509614757702, 8.691521330119459, 0, 15.84431199687321, 18.89162482598436, 22.4039850392055),
    @(10.02966061733435, 5.542578669922527, 11.25460463087819, 0, 30.94820226971435, 29.84869459387691, 14.56257616020193, 0, 11.37694575505755, 9.001494437986599, 0, 16.00528319733916, 18.80601438653344, 22.31390409642875),
    @(10.22888130343982, 5.615390478016069, 11.29969960919308, 0, 30.93753669617863, 29.8303843733115, 14.53731812062769, 0, 11.37039259288481, 9.139880347123942, 0, 16.08045192533272, 18.76871480719017, 22.27725036075044),
    @(10.30345714093136, 5.642651178653241, 11.31701425494616, 0, 30.93450211873496, 29.8246948616803, 14.52805366029908, 0, 11.36815563215169, 9.19185846615064, 0, 16.10917716509602, 18.75482582359886, 22.26399318747038),
    @(10.28743539208579, 5.636794215357471, 11.3132748246081, 0, 30.93511102038628, 29.82586483326472, 14.53003557819058, 0, 11.36862653075305, 9.180683714991385, 0, 16.10297938962241, 18.75780660398922, 22.26682064347608),
    @(10.23503440633902, 5.617639570056386, 11.30111938734455, 0, 30.93726691841229, 29.82989134709789, 14.5365499121466, 0, 11.37020365747792, 9.144165414431338, 0, 16.08281005689842, 18.76756743689649, 22.27614719813235),
    @(10.20282273267617, 5.60586573582262, 11.29370451056046, 0, 30.93871822372438, 29.83251979221852, 14.54057922135231, 0, 11.37120153206641, 9.121740089205973, 0, 16.0704891148598, 18.77357687604729, 22.28194111592446),
    @(10.01652370167615, 5.53777755453815, 11.25169150693329, 0, 30.94903986372527, 29.8500651944312, 14.56426885147252, 0, 11.37740826748608, 8.992393600051843, 0, 16.00040822679659, 18.80848501866759, 22.3163866066627),
    @(9.900767313740582, 5.495470664892718, 11.22635405865842, 0, 30.95716134185573, 29.86304231307588, 14.57933655938166, 0, 11.3816519600545, 8.9123364787488, 0, 15.95790005132991, 18.83032065444807, 22.33862610519606),
    @(9.833670560102973, 5.470945065567119, 11.21194446245763, 0, 30.96249039201083, 29.87131904101686, 14.58819972786349, 0, 11.38425319226677, 8.866045598109176, 0, 15.93363442815362, 18.84303483234289, 22.35182467038623),
    @(9.81086623890652, 5.462608717077947, 11.20709409097983, 0, 30.96440770217553, 29.87426087966305, 14.59123441931052, 0, 11.38516147974973, 8.850331902599782, 0, 15.92545066564112, 18.84736626892296, 22.35636334972688),
    @(9.913143786268078, 5.499994282960022, 11.22903439621002, 0, 30.95622872032517, 29.86157675146319, 14.57771222774886, 0, 11.38118361640282, 8.920884325779117, 0, 15.96240621542107, 18.82798018966336, 22.33621654326863),
    @(10.25044983403513, 5.623274339660179, 11.30468336063652, 0, 30.93660643137943, 29.82867487977078, 14.53462834799583, 0, 11.36973378270234, 9.154903645941639, 0, 16.08872736142812, 18.76469405858969, 22.27339085316854),
    @(10.46583122314476, 5.702021895699424, 11.35550709249551, 0, 30.92963497700007, 29.81442419368351, 14.50822032769972, 0, 11.36367602817076, 9.305347798321584, 0, 16.172794791162, 18.72470546882947, 22.23596118075524),
    @(10.35136300654654, 5.660165008764576, 11.32825874445778, 0, 30.93282057181001, 29.82136581388931, 14.52215473281228, 0, 11.36677888789025, 9.225296933455205, 0, 16.12779463365132, 18.74592286934174, 22.25560561348216),
    @(9.907550081032547, 5.497949785628894, 11.2278221249851, 0, 30.95664830293608, 29.86223679026556, 14.57844596338898, 0, 11.38139485155646, 8.917020665747684, 0, 15.96036843866828, 18.82903781351558, 22.33730462050168),
    @(9.409415489658009, 5.315720508366059, 11.12578943403433, 0, 31.00749624145874, 29.93737371024552, 14.6466831570744, 0, 11.40329286815645, 8.575385849012971, 0, 15.78687646770586, 19.00013952161483, 22.44098607698325)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item(2 + $r, 2 + $c).Value = $row[$c]
    }
}
